# "Generate Report for Handback"
# The 88a16c70-8afc-41f5-81ef-23758a2bce06 localization entry has now been
# handed back (it is in sync with en-US). Reflect that on the Overview
# sheet and on each per-language detail sheet: flip the Status column,
# populate the Latest Target File / Latest Handback File columns (with
# hyperlinks, like the other file-name columns already have), and stamp
# the Latest Handback DateTime.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus

# ---- Per-language detail sheets --------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; TargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/c7e45d08b2086fcf838eaa4e54302e31769bcba8/e2e/88a16c70-8afc-41f5-81ef-23758a2bce06.md"; TargetName = "88a16c70-8afc-41f5-81ef-23758a2bce06.md"; HandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cd35831e7c23e991c963f14291c922b245b7f36/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/88a16c70-8afc-41f5-81ef-23758a2bce06.1d84dd4d5cdc3eb60110f8c032bb8103ede39b42.zh-cn.xlf"; HandbackName = "88a16c70-8afc-41f5-81ef-23758a2bce06.1d84dd4d5cdc3eb60110f8c032bb8103ede39b42.zh-cn.xlf"; HandbackTime = "2016-03-14 00:50:16" },
    @{ Sheet = "de-de"; TargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/c7e45d08b2086fcf838eaa4e54302e31769bcba8/e2e/88a16c70-8afc-41f5-81ef-23758a2bce06.md"; TargetName = "88a16c70-8afc-41f5-81ef-23758a2bce06.md"; HandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4382e39a8a4ea7b0cbcec1e591257c8d090c1ad7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/88a16c70-8afc-41f5-81ef-23758a2bce06.1d84dd4d5cdc3eb60110f8c032bb8103ede39b42.de-de.xlf"; HandbackName = "88a16c70-8afc-41f5-81ef-23758a2bce06.1d84dd4d5cdc3eb60110f8c032bb8103ede39b42.de-de.xlf"; HandbackTime = "2016-03-14 00:50:22" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Remember the existing (row 3 / "d0dd5a4b...") hyperlinks so they can
    # be re-created after the new row-2 links, keeping everything in a
    # single, row-ordered Hyperlinks collection.
    $existing = @()
    foreach ($h in $ws.Hyperlinks) {
        $existing += , @($h.Range.Address($false, $false), $h.Address, $h.TextToDisplay)
    }
    $ws.Hyperlinks.Delete()

    foreach ($item in $existing) {
        $addr = $item[0]
        if ($addr -eq "A2" -or $addr -eq "B2" -or $addr -eq "D2") {
            $ws.Hyperlinks.Add($ws.Range($addr), $item[1], "", "", $item[2])
        }
    }

    # New columns for the handed-back file: Latest Target File (F) and
    # Latest Handback File (G).
    $ws.Range("F2").Value = $lang.TargetName
    $ws.Hyperlinks.Add($ws.Range("F2"), $lang.TargetUrl, "", "", $lang.TargetName)

    $ws.Range("G2").Value = $lang.HandbackName
    $ws.Hyperlinks.Add($ws.Range("G2"), $lang.HandbackUrl, "", "", $lang.HandbackName)

    foreach ($item in $existing) {
        $addr = $item[0]
        if ($addr -eq "A3" -or $addr -eq "B3" -or $addr -eq "D3") {
            $ws.Hyperlinks.Add($ws.Range($addr), $item[1], "", "", $item[2])
        }
    }

    # Status + Latest Handback DateTime for the handed-back row.
    $ws.Range("C2").Value = $newStatus
    $ws.Range("H2").Value = $lang.HandbackTime
}
